$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price and 1h volume change) per diff
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.723.40"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.722.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.31"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.95"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.721.43"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.00%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.05"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.346.27"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.719.70"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.714.91"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.27"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.16"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "493.34"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.56"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +15.18%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.07"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.45"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.93%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.46"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.868.99"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.657.95"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.324"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "437.40"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.85"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.47"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.58"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.21"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.767.72"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.82%  "
